$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") '26.854.94'
Set-TextValue $ws.Range("E2") '  +1.49%  '

Set-TextValue $ws.Range("D3") '1.730.42'
Set-TextValue $ws.Range("E3") '  +0.32%  '

Set-TextValue $ws.Range("D4") '0.9990'
Set-TextValue $ws.Range("E4") '  +0.24%  '

Set-TextValue $ws.Range("D5") '240.76'
Set-TextValue $ws.Range("E5") '  -0.87%  '

Set-TextValue $ws.Range("D6") '0.9994'
Set-TextValue $ws.Range("E6") '  +0.21%  '

Set-TextValue $ws.Range("D7") '0.4832'
Set-TextValue $ws.Range("E7") '  -1.14%  '

Set-TextValue $ws.Range("E8") '  -0.52%  '

Set-TextValue $ws.Range("D9") '0.06182'
Set-TextValue $ws.Range("E9") '  -0.25%  '

Set-TextValue $ws.Range("D10") '1.729.19'
Set-TextValue $ws.Range("E10") '  +0.38%  '

Set-TextValue $ws.Range("E11") '  +2.14%  '

Set-TextValue $ws.Range("D12") '0.06876'
Set-TextValue $ws.Range("E12") '  -1.57%  '

Set-TextValue $ws.Range("D13") '0.6044'
Set-TextValue $ws.Range("E13") '  +0.51%  '

Set-TextValue $ws.Range("D14") '4.461'
Set-TextValue $ws.Range("E14") '  -1.32%  '

Set-TextValue $ws.Range("D15") '77.02'
Set-TextValue $ws.Range("E15") '  -0.28%  '

Set-TextValue $ws.Range("D16") '0.9996'
Set-TextValue $ws.Range("E16") '  +0.23%  '

Set-TextValue $ws.Range("D17") '26.834.07'
Set-TextValue $ws.Range("E17") '  +1.48%  '

Set-TextValue $ws.Range("D18") '0.9992'
Set-TextValue $ws.Range("E18") '  +0.24%  '

Set-TextValue $ws.Range("D19") '0.000007133'
Set-TextValue $ws.Range("E19") '  -0.42%  '

Set-TextValue $ws.Range("D20") '11.39'
Set-TextValue $ws.Range("E20") '  +0.61%  '

Set-TextValue $ws.Range("D21") '1.951.07'
Set-TextValue $ws.Range("E21") '  +0.41%  '

Set-TextValue $ws.Range("D22") '4.414'
Set-TextValue $ws.Range("E22") '  -1.11%  '

Set-TextValue $ws.Range("D23") '8.460'
Set-TextValue $ws.Range("E23") '  -0.69%  '

Set-TextValue $ws.Range("D24") '5.072'
Set-TextValue $ws.Range("E24") '  -0.65%  '

Set-TextValue $ws.Range("D25") '140.43'
Set-TextValue $ws.Range("E25") '  +2.04%  '

Set-TextValue $ws.Range("D26") '15.22'
Set-TextValue $ws.Range("E26") '  -0.37%  '

Set-TextValue $ws.Range("D27") '1.793'
Set-TextValue $ws.Range("E27") '  +2.85%  '

Set-TextValue $ws.Range("D28") '106.49'
Set-TextValue $ws.Range("E28") '  -0.11%  '

Set-TextValue $ws.Range("D29") '1.374'
Set-TextValue $ws.Range("E29") '  -2.48%  '

Set-TextValue $ws.Range("D30") '3.955'
Set-TextValue $ws.Range("E30") '  +1.07%  '

Set-TextValue $ws.Range("D31") '0.07920'
Set-TextValue $ws.Range("E31") '  -1.26%  '

Set-TextValue $ws.Range("D32") '3.664'
Set-TextValue $ws.Range("E32") '  +0.47%  '

Set-TextValue $ws.Range("D33") '0.04551'
Set-TextValue $ws.Range("E33") '  +1.30%  '

Set-TextValue $ws.Range("D34") '2.596'
Set-TextValue $ws.Range("E34") '  -0.03%  '

Set-TextValue $ws.Range("D35") '1.002'
Set-TextValue $ws.Range("E35") '  +0.15%  '

Set-TextValue $ws.Range("D36") '0.6169'
Set-TextValue $ws.Range("E36") '  -1.27%  '

Set-TextValue $ws.Range("D37") '0.9244'
Set-TextValue $ws.Range("E37") '  -0.40%  '

$ws.Range("B38").Value = 'MXToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue $ws.Range("D38") '2.446'
Set-TextValue $ws.Range("E38") '  +2.54%  '

$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws.Range("D39") '1.996'
Set-TextValue $ws.Range("E39") '  +0.54%  '

Set-TextValue $ws.Range("D40") '0.9991'
Set-TextValue $ws.Range("E40") '  +0.25%  '

Set-TextValue $ws.Range("D41") '0.01493'
Set-TextValue $ws.Range("E41") '  +0.87%  '

Set-TextValue $ws.Range("D42") '5.638'
Set-TextValue $ws.Range("E42") '  +4.86%  '

Set-TextValue $ws.Range("D43") '99.73'
Set-TextValue $ws.Range("E43") '  -0.08%  '

Set-TextValue $ws.Range("D44") '0.3830'
Set-TextValue $ws.Range("E44") '  -0.35%  '

Set-TextValue $ws.Range("D45") '6.779'
Set-TextValue $ws.Range("E45") '  -1.91%  '

Set-TextValue $ws.Range("D46") '0.1157'
Set-TextValue $ws.Range("E46") '  -0.49%  '

Set-TextValue $ws.Range("E47") '  -0.04%  '

Set-TextValue $ws.Range("D48") '7.897'
Set-TextValue $ws.Range("E48") '  +2.40%  '

Set-TextValue $ws.Range("D49") '30.07'
Set-TextValue $ws.Range("E49") '  -1.07%  '

Set-TextValue $ws.Range("D50") '1.242'
Set-TextValue $ws.Range("E50") '  +1.46%  '

Set-TextValue $ws.Range("D51") '51.58'
Set-TextValue $ws.Range("E51") '  +0.68%  '

Write-Output "Applied cryptos update"